# Apply edit: add a new "Slacks" product row to the Products sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")

# New row 5 data: name=Slacks, price=60, category=Pants (column C/imageUrl, E left blank)
$ws.Range("A5").Value = "Slacks"
$ws.Range("B5").Value = 60
$ws.Range("D5").Value = "Pants"

# Update the active selection to match the saved workbook state
$ws.Activate()
$ws.Range("A6").Select()
